$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A3").Value = -22.03070000000002
$ws.Range("D3").Value = -7.396499999999995
$ws.Range("A21").Value = -20.08829999999999
$ws.Range("A23").Value = -20.28199999999998
$ws.Range("D24").Value = -7.586899999999997
$ws.Range("A25").Value = -21.70049999999998
$ws.Range("B27").Value = 6.024600000000003
$ws.Range("B31").Value = 5.869800000000001
$ws.Range("B39").Value = 9.919699999999999
$ws.Range("B48").Value = 5.151000000000005
$ws.Range("B51").Value = 5.431799999999998
$ws.Range("B52").Value = 5.189800000000002
$ws.Range("A53").Value = -21.84640000000001
$ws.Range("B55").Value = 6.008599999999997
$ws.Range("B56").Value = 5.192699999999999
$ws.Range("A57").Value = -22.3569
$ws.Range("B57").Value = 4.825999999999997
$ws.Range("D57").Value = -8.6332
$ws.Range("A59").Value = -22.3141
$ws.Range("D61").Value = -7.643599999999997
$ws.Range("A69").Value = -21.65269999999999
$ws.Range("D70").Value = -7.225999999999996
$ws.Range("B73").Value = 8.546999999999999
$ws.Range("A79").Value = -20.51150000000001
$ws.Range("A83").Value = -21.8828
$ws.Range("D86").Value = -7.573299999999992
$ws.Range("B89").Value = 5.195599999999996
$ws.Range("B90").Value = 5.884000000000004
$ws.Range("A93").Value = -21.2361
$ws.Range("D98").Value = -8.513799999999998
$ws.Range("D100").Value = -8.296799999999999
$ws.Range("D102").Value = -7.719399999999996
